$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the paragraph "Remember that this is a public repository
#    - your changes could be seen by anyone who looks!" with two runs
#    of (CJK-tagged) text: "yeza" and "12", keeping the paragraph's
#    border formatting (pBdr) intact.
# ------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("Remember that this is a public repository", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $findRng.Paragraphs(1).Range

    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
        'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
        'w14:paraId="2570EC4D" w14:textId="0270A0CC" w:rsidR="00091943" ' + `
        'w:rsidRDefault="00091943" w:rsidP="00EF6287">' + `
        '<w:pPr>' + `
            '<w:pBdr>' + `
                '<w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>' + `
                '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' + `
                '<w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>' + `
                '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' + `
            '</w:pBdr>' + `
        '</w:pPr>' + `
        '<w:r>' + `
            '<w:rPr>' + `
                '<w:rFonts w:hint="eastAsia"/>' + `
                '<w:lang w:eastAsia="zh-CN"/>' + `
            '</w:rPr>' + `
            '<w:t>yeza</w:t>' + `
        '</w:r>' + `
        '<w:r>' + `
            '<w:rPr>' + `
                '<w:lang w:eastAsia="zh-CN"/>' + `
            '</w:rPr>' + `
            '<w:t>12</w:t>' + `
        '</w:r>' + `
    '</w:p>'

    $target.InsertXML($newParaXml)
}

# ------------------------------------------------------------------
# 2) Remove the three paragraphs that followed it in full:
#      "Add some comments about Version management outside this
#       border, or just add some text so there is a change to this
#       file."
#      "Remember that your GitHub user ID must be submitted in your
#       assignment report!"
#      "Once you've changed follow the next step in your assignment
#       task."
# ------------------------------------------------------------------
$startRng = $d.Content
$startFound = $startRng.Find.Execute("Add some comments about Version management", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRng = $d.Content
$endFound = $endRng.Find.Execute("Once you", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($startFound -and $endFound) {
    $startPara = $startRng.Paragraphs(1)
    $endPara = $endRng.Paragraphs(1)
    $delRng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRng.Delete()
}
